$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.241.37'
$ws.Range('E2').Value = '  -0.69%  '
$ws.Range('D3').Value = '2.270.67'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '308.81'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.92'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  -0.94%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.07'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0804'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.09%  '
$ws.Range('E12').Value = '  +0.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.83'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').Value = '2.622.37'
$ws.Range('E14').Value = '  -0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.66'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.53%  '
$ws.Range('D16').Value = '2.267.91'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('D18').Value = '42.103.10'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.21'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -4.28%  '
$ws.Range('D20').Value = '0.0₃0903'
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('E21').Value = '  -0.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.62'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.31'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.86%  '
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.53'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.21'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.55'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '163.54'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +1.49%  '
$ws.Range('E32').Value = '  -1.46%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.63'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.54%  '
$ws.Range('E36').Value = '  -2.60%  '
$ws.Range('E37').Value = '  -0.36%  '
$ws.Range('E38').Value = '  -3.80%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.115'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.07%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.82'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.19'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('E42').Value = '  -6.47%  '
$ws.Range('D43').Value = '1.947.88'
$ws.Range('E44').Value = '  -2.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.76'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.95'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.74%  '
$ws.Range('E47').Value = '  -4.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.75'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').Value = '2.494.72'
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '92.18'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('E51').Value = '  -1.88%  '
